$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.322.53"
$ws.Range("E2").Value = "  +0.24%  "

$ws.Range("D3").Value = "'1.874.99"
$ws.Range("E3").Value = "  +0.24%  "

$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'0.7121"
$ws.Range("E5").Value = "  +0.13%  "

$ws.Range("D6").Value = "'242.14"
$ws.Range("E6").Value = "  +0.28%  "

$ws.Range("D7").Value = "'0.9997"
$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("D8").Value = "'0.07876"
$ws.Range("E8").Value = "  +2.38%  "

$ws.Range("D9").Value = "'0.3118"
$ws.Range("E9").Value = "  +0.40%  "

$ws.Range("D10").Value = "'25.25"
$ws.Range("E10").Value = "  +0.93%  "

$ws.Range("E11").Value = "  +0.53%  "

$ws.Range("D12").Value = "'1.885.24"
$ws.Range("E12").Value = "  +0.59%  "

$ws.Range("D13").Value = "'5.245"
$ws.Range("E13").Value = "  +0.59%  "

$ws.Range("D14").Value = "'0.7188"
$ws.Range("E14").Value = "  +1.28%  "

$ws.Range("D15").Value = "'91.33"
$ws.Range("E15").Value = "  +0.18%  "

$ws.Range("D16").Value = "'6.148"
$ws.Range("E16").Value = "  +3.64%  "

$ws.Range("D17").Value = "'0.000008361"
$ws.Range("E17").Value = "  +1.21%  "

$ws.Range("D18").Value = "'29.315.23"
$ws.Range("E18").Value = "  +0.09%  "

$ws.Range("D19").Value = "'240.90"
$ws.Range("E19").Value = "  -0.43%  "

$ws.Range("D20").Value = "'13.23"
$ws.Range("E20").Value = "  +0.54%  "

$ws.Range("D21").Value = "'2.121.39"
$ws.Range("E21").Value = "  -0.53%  "

$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  +0.15%  "

$ws.Range("D23").Value = "'7.787"
$ws.Range("E23").Value = "  -0.62%  "

$ws.Range("E24").Value = "  -0.05%  "

$ws.Range("D25").Value = "'0.1602"
$ws.Range("E25").Value = "  -1.27%  "

$ws.Range("D26").Value = "'162.88"
$ws.Range("E26").Value = "  -0.20%  "

$ws.Range("D27").Value = "'9.058"
$ws.Range("E27").Value = "  +0.68%  "

$ws.Range("D28").Value = "'18.54"
$ws.Range("E28").Value = "  +0.24%  "

$ws.Range("E29").Value = "  +0.33%  "

$ws.Range("D30").Value = "'4.421"
$ws.Range("E30").Value = "  +0.25%  "

$ws.Range("D31").Value = "'4.348"
$ws.Range("E31").Value = "  +0.22%  "

$ws.Range("D32").Value = "'1.230"
$ws.Range("E32").Value = "  -3.77%  "

$ws.Range("D33").Value = "'0.05366"
$ws.Range("E33").Value = "  +2.44%  "

$ws.Range("D34").Value = "'1.949"
$ws.Range("E34").Value = "  +1.39%  "

$ws.Range("D35").Value = "'1.177"
$ws.Range("E35").Value = "  +0.53%  "

$ws.Range("D36").Value = "'0.7473"
$ws.Range("E36").Value = "  -0.26%  "

$ws.Range("D37").Value = "'2.687"
$ws.Range("E37").Value = "  +0.22%  "

$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "'1.295.48"
$ws.Range("E38").Value = "  +12.41%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.01877"
$ws.Range("E39").Value = "  +1.08%  "

$ws.Range("E40").Value = "  +0.84%  "

$ws.Range("D41").Value = "'6.516"
$ws.Range("E41").Value = "  +2.41%  "

$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "'110.32"
$ws.Range("E42").Value = "  +5.51%  "

$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'0.8918"
$ws.Range("E43").Value = "  +0.65%  "

$ws.Range("D44").Value = "'72.70"
$ws.Range("E44").Value = "  -0.38%  "

$ws.Range("D45").Value = "'0.00000000134"
$ws.Range("E45").Value = "  +14.54%  "

$ws.Range("D46").Value = "'0.9997"
$ws.Range("E46").Value = "  -0.04%  "

$ws.Range("D47").Value = "'2.016.61"
$ws.Range("E47").Value = "  -0.56%  "

$ws.Range("D48").Value = "'1.803"
$ws.Range("E48").Value = "  +0.70%  "

$ws.Range("D49").Value = "'0.5188"
$ws.Range("E49").Value = "  +0.09%  "

$ws.Range("D50").Value = "'9.465"
$ws.Range("E50").Value = "  +1.14%  "

$ws.Range("D51").Value = "'0.4362"
$ws.Range("E51").Value = "  +1.71%  "
